$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (|S*|/n) over the 10 data rows, shown in bold.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Summary rows 14-17: headline metrics in column A (new shared strings) with
# their values (bold, 12pt, vertically centered) in column B.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$summary = $ws.Range("B14:B17")
$summary.Font.Bold = $true
$summary.Font.Size = 12
$summary.VerticalAlignment = -4108

$ws.Rows(14).RowHeight = 15.6
$ws.Rows(15).RowHeight = 15.6
$ws.Rows(16).RowHeight = 15.6
$ws.Rows(17).RowHeight = 15.6

# Page setup recorded by Excel the last time the sheet was printed/saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection on J12, matching the last-active cell recorded by Excel.
$ws.Range("J12").Select()

